$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 800.625
$ws.Range("I98").Value = 800.625
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 800.625
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 697.375
$ws.Range("N98").Value = $null

# Row 122
$ws.Range("H122").Value = 800.625
$ws.Range("I122").Value = 800.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2401.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 48.125
$ws.Range("N122").Value = $null

# Row 129
$ws.Range("H129").Value = 864.3946999999999
$ws.Range("J129").Value = 874.2432
$ws.Range("L129").Value = 2622.7296
$ws.Range("N129").Value = -12622.7296

$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 6684.5
$ws.Range("J44").Value = 6684.5
$ws.Range("L44").Value = 6684.5
$ws.Range("N44").Value = -7660.5

# Row 55
$ws.Range("H55").Value = 10348
$ws.Range("J55").Value = 10685
$ws.Range("L55").Value = 10685
$ws.Range("N55").Value = -11315

# Row 61
$ws.Range("H61").Value = 1696.8572
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 1813
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 1813
$ws.Range("M61").Value = -788
$ws.Range("N61").Value = -2237

# Row 63
$ws.Range("H63").Value = 2664.8
$ws.Range("I63").Value = 2349.7144
$ws.Range("K63").Value = 2349.7144
$ws.Range("M63").Value = -1663.7144

# Row 66
$ws.Range("H66").Value = 2664.8
$ws.Range("I66").Value = 2349.7144
$ws.Range("K66").Value = 11748.572
$ws.Range("M66").Value = -8316.572

# Row 74
$ws.Range("H74").Value = 1953.1082
$ws.Range("I74").Value = 1076
$ws.Range("K74").Value = 1076
$ws.Range("M74").Value = -202

# Row 77
$ws.Range("H77").Value = 1953.1082
$ws.Range("I77").Value = 1076
$ws.Range("K77").Value = 5380
$ws.Range("M77").Value = -1012

# Row 86
$ws.Range("H86").Value = 56542.668
$ws.Range("J86").Value = 56542.668
$ws.Range("L86").Value = 56542.668
$ws.Range("N86").Value = -58914.668

# Row 89
$ws.Range("H89").Value = 56542.668
$ws.Range("J89").Value = 56542.668
$ws.Range("L89").Value = 169628.004
$ws.Range("N89").Value = -181484.004

# Row 110
$ws.Range("H110").Value = 47672150
$ws.Range("I110").Value = 50055708
$ws.Range("K110").Value = 50055708
$ws.Range("M110").Value = -50053663

# Row 136
$ws.Range("H136").Value = 1696.8572
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1813
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 5439
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -10539

$ws = $wb.Worksheets.Item("BSM")
# Row 23
$ws.Range("H23").Value = 4014
$ws.Range("J23").Value = 4014
$ws.Range("L23").Value = 4014
$ws.Range("N23").Value = -4580

# Row 35
$ws.Range("H35").Value = 17537
$ws.Range("J35").Value = 17537
$ws.Range("L35").Value = 17537
$ws.Range("N35").Value = -18157

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 28433.79
$ws.Range("I31").Value = 37060.18
$ws.Range("J31").Value = 4279.9
$ws.Range("K31").Value = 37060.18
$ws.Range("L31").Value = 4279.9
$ws.Range("M31").Value = -36765.18
$ws.Range("N31").Value = -4869.9

# Row 34
$ws.Range("H34").Value = 28433.79
$ws.Range("I34").Value = 37060.18
$ws.Range("J34").Value = 4279.9
$ws.Range("K34").Value = 37060.18
$ws.Range("L34").Value = 4279.9
$ws.Range("M34").Value = -36858.18
$ws.Range("N34").Value = -4683.9

# Row 99
$ws.Range("H99").Value = 27748.25
$ws.Range("J99").Value = 35337.668
$ws.Range("L99").Value = 35337.668
$ws.Range("N99").Value = -38333.668

# Row 126
$ws.Range("H126").Value = 27748.25
$ws.Range("J126").Value = 35337.668
$ws.Range("L126").Value = 106013.004
$ws.Range("N126").Value = -110953.004

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 16970.23
$ws.Range("I5").Value = 3102
$ws.Range("J5").Value = 19491.727
$ws.Range("K5").Value = 9306
$ws.Range("L5").Value = 58475.181
$ws.Range("M5").Value = -9194
$ws.Range("N5").Value = -58699.181

# Row 34
$ws.Range("H34").Value = 2645.4285
$ws.Range("J34").Value = 2837.3845
$ws.Range("L34").Value = 8512.1535
$ws.Range("N34").Value = -8680.1535

# Row 55
$ws.Range("H55").Value = 6761.95
$ws.Range("I55").Value = 650
$ws.Range("J55").Value = 7441.0557
$ws.Range("K55").Value = 1950
$ws.Range("L55").Value = 22323.1671
$ws.Range("M55").Value = -1773
$ws.Range("N55").Value = -22677.1671

# Row 135
$ws.Range("H135").Value = 16970.23
$ws.Range("I135").Value = 3102
$ws.Range("J135").Value = 19491.727
$ws.Range("K135").Value = 27918
$ws.Range("L135").Value = 175425.543
$ws.Range("M135").Value = -25383
$ws.Range("N135").Value = -180495.543

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 5300

# Row 71
$ws.Range("H71").Value = 5300

# Row 74
$ws.Range("H74").Value = 15444
$ws.Range("J74").Value = 16124.5
$ws.Range("L74").Value = 16124.5
$ws.Range("N74").Value = -18120.5

# Row 77
$ws.Range("H77").Value = 15444
$ws.Range("J77").Value = 16124.5
$ws.Range("L77").Value = 48373.5
$ws.Range("N77").Value = -58357.5

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872

# Row 78
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360

# Row 132
$ws.Range("H132").Value = 1900.3334
$ws.Range("I132").Value = 1130.3
$ws.Range("J132").Value = 2862.875
$ws.Range("K132").Value = 3390.9
$ws.Range("L132").Value = 8588.625
$ws.Range("M132").Value = -860.8999999999996
$ws.Range("N132").Value = -13648.625

Write-Output "Applied all profit-sheet updates"